$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "27.616.87", "1.001") are not coerced into numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.616.87'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '1.802.25'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.52%  '

$ws.Range("D5").Value = '337.01'
$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("D6").Value = '0.9976'
$ws.Range("E6").Value = '  -0.57%  '

$ws.Range("D7").Value = '0.3935'
$ws.Range("E7").Value = '  +3.85%  '

$ws.Range("D8").Value = '0.3467'
$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("D9").Value = '48.19'
$ws.Range("E9").Value = '  -2.25%  '

$ws.Range("D10").Value = '1.201'
$ws.Range("E10").Value = '  -0.77%  '

$ws.Range("D11").Value = '0.07514'
$ws.Range("E11").Value = '  -0.64%  '

$ws.Range("D12").Value = '0.9981'
$ws.Range("E12").Value = '  -0.54%  '

$ws.Range("D13").Value = '22.18'
$ws.Range("E13").Value = '  +1.09%  '

$ws.Range("D14").Value = '6.521'
$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("D15").Value = '1.800.85'
$ws.Range("E15").Value = '  +0.23%  '

$ws.Range("D16").Value = '7.159'
$ws.Range("E16").Value = '  +1.52%  '

$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("D18").Value = '0.06689'
$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("D19").Value = '84.83'
$ws.Range("E19").Value = '  -0.05%  '

$ws.Range("D20").Value = '0.9971'
$ws.Range("E20").Value = '  -0.58%  '

$ws.Range("D21").Value = '17.77'
$ws.Range("E21").Value = '  +1.81%  '

$ws.Range("D22").Value = '6.559'
$ws.Range("E22").Value = '  +1.09%  '

$ws.Range("D23").Value = '27.626.29'
$ws.Range("E23").Value = '  +0.79%  '

$ws.Range("D24").Value = '12.82'
$ws.Range("E24").Value = '  +2.06%  '

$ws.Range("D25").Value = '2.392'
$ws.Range("E25").Value = '  -2.51%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = '1.477'
$ws.Range("E26").Value = '  -0.98%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.522'
$ws.Range("E27").Value = '  -2.25%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '21.26'
$ws.Range("E28").Value = '  -0.96%  '

$ws.Range("D29").Value = '155.70'
$ws.Range("E29").Value = '  +3.34%  '

$ws.Range("D30").Value = '2.000.80'
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").Value = '135.48'
$ws.Range("E31").Value = '  +1.35%  '

$ws.Range("D32").Value = '4.028'
$ws.Range("E32").Value = '  -1.09%  '

$ws.Range("D33").Value = '6.032'
$ws.Range("E33").Value = '  -1.39%  '

$ws.Range("D34").Value = '0.08804'
$ws.Range("E34").Value = '  +1.39%  '

$ws.Range("D35").Value = '13.14'
$ws.Range("E35").Value = '  -0.90%  '

$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = '5.474'
$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02424'
$ws.Range("E37").Value = '  +2.63%  '

$ws.Range("D38").Value = '1.613'
$ws.Range("E38").Value = '  -4.04%  '

$ws.Range("D39").Value = '0.06498'
$ws.Range("E39").Value = '  +2.23%  '

$ws.Range("D40").Value = '0.6852'
$ws.Range("E40").Value = '  -0.41%  '

$ws.Range("D41").Value = '0.2220'
$ws.Range("E41").Value = '  +0.37%  '

$ws.Range("D42").Value = '1.259'
$ws.Range("E42").Value = '  -1.40%  '

$ws.Range("D43").Value = '8.463'
$ws.Range("E43").Value = '  -4.05%  '

$ws.Range("D44").Value = '14.65'
$ws.Range("E44").Value = '  +0.91%  '

$ws.Range("D45").Value = '0.6470'
$ws.Range("E45").Value = '  +0.19%  '

$ws.Range("D46").Value = '0.9969'
$ws.Range("E46").Value = '  -0.56%  '

$ws.Range("D47").Value = '3.854'
$ws.Range("E47").Value = '  -0.06%  '

$ws.Range("D48").Value = '2.155'
$ws.Range("E48").Value = '  +0.96%  '

$ws.Range("D49").Value = '132.56'
$ws.Range("E49").Value = '  +1.24%  '

$ws.Range("D50").Value = '0.07217'
$ws.Range("E50").Value = '  +0.13%  '

$ws.Range("D51").Value = '80.11'
$ws.Range("E51").Value = '  +0.56%  '

# Restore default cell style on column D (NumberFormat change above
# would otherwise leave a lingering explicit style on the cells).
$ws.Range("D2:D51").Style = "Normal"
